$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$data = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# Sheet "Data" - new row for 2021 (label first, to match shared-string order)
# ---------------------------------------------------------------------------
$data.Range("A59").Value = "2021............................................................................. ."

# ---------------------------------------------------------------------------
# Sheet "About" - point the annual CPI reference at the newer BLS PDF
# ---------------------------------------------------------------------------
foreach ($h in $about.Hyperlinks) {
    $h.Delete()
}
$about.Range("B6").Value = "https://www.bls.gov/cpi/tables/supplemental-files/historical-cpi-u-202212.pdf"
$about.Range("B7").Value = "Pages 4 and 5"

# ---------------------------------------------------------------------------
# Sheet "Data" - finish row 59 (2021) and add row 60 (2022)
# ---------------------------------------------------------------------------
$data.Range("B59").Value = 266.23599999999999
$data.Range("C59").Value = 275.70299999999997
$data.Range("D59").Value = 270.97000000000003
$data.Range("E59").Value = 7
$data.Range("F59").Value = 4.7
$data.Range("G59").Formula = "=`$D`$50/D59"
$data.Range("G59").Style = $data.Range("G58").Style
$data.Range("G59").NumberFormat = $data.Range("G58").NumberFormat

$data.Range("A60").Value = "2022............................................................................. ."
$data.Range("B60").Value = 288.34699999999998
$data.Range("C60").Value = 296.96300000000002
$data.Range("D60").Value = 292.65499999999997
$data.Range("E60").Value = 6.5
$data.Range("F60").Value = 8
$data.Range("G60").Formula = "=`$D`$50/D60"
$data.Range("G60").Style = $data.Range("G58").Style
$data.Range("G60").NumberFormat = $data.Range("G58").NumberFormat

# Blank row 61, then a lone styled (but empty) cell at G62
$data.Range("G62").Style = $data.Range("G58").Style
$data.Range("G62").NumberFormat = $data.Range("G58").NumberFormat
$data.Range("G62").ClearContents()

# ---------------------------------------------------------------------------
# Window / selection state
# ---------------------------------------------------------------------------
$data.Activate()
$data.Range("B63").Select()

$about.Activate()
$about.Range("B8").Select()
